$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows above row 7, shifting the existing data (and all
# subsequent rows) down by two rows.
$ws.Rows("7:8").Insert()

# Populate the two new rows with the new "input"/"output" configuration_fxe
# entries for the FxE matrix. Shared strings are created in the order the
# values are first written, so write column C (Parameter) for both rows
# before column D (Type) etc. to match the expected string table order.
$ws.Cells.Item(7, 1).Value = "CHE"
$ws.Cells.Item(7, 2).Value = "conv_chp_coal"
$ws.Cells.Item(7, 3).Value = "input"

$ws.Cells.Item(8, 1).Value = "CHE"
$ws.Cells.Item(8, 2).Value = "conv_chp_coal"
$ws.Cells.Item(8, 3).Value = "output"

$ws.Cells.Item(7, 4).Value = "configuration_fxe"
$ws.Cells.Item(8, 4).Value = "configuration_fxe"

$ws.Cells.Item(7, 6).Value = "coal"
$ws.Cells.Item(8, 6).Value = "elecsupply"

$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(8, 7).Value = 1

# Re-apply the worksheet AutoFilter over the new, larger range (toggle off
# first since the range is already filtered and simply re-selecting would
# switch it off instead of updating the range).
$ws.AutoFilterMode = $false
$ws.Range("A5:L852").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$5:`$L`$852"

# Match the author's final selection state.
$ws.Range("F7:F8").Select()
